# Updated cryptos list with GitHub Actions — refreshed Price / Volume(1h)
# figures, plus the Arweave <-> EnergySwap row-identity swap (rows 45/46).
#
# Some Price values (column D) are plain decimal numbers (e.g. "1.00",
# "602.83") that Excel would normally auto-convert to numeric values,
# silently dropping trailing zeros or introducing floating-point noise
# (e.g. "602.83" -> 602.83000000000004). To keep these as literal text
# exactly as they appear in the source data, cells flagged `Quote = $true`
# are written with a leading apostrophe (Excel's "treat as text" marker,
# which is not itself stored as part of the value) and then have their
# style reset back to "Normal" so no stray text-format style lingers on
# the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '67.801.29'; Quote = $false }
    @{ Cell = "E2"; Value = '  +0.28%  '; Quote = $false }
    @{ Cell = "D3"; Value = '3.826.13'; Quote = $false }
    @{ Cell = "E3"; Value = '  +1.21%  '; Quote = $false }
    @{ Cell = "E4"; Value = '  +0.12%  '; Quote = $false }
    @{ Cell = "D5"; Value = '602.83'; Quote = $true }
    @{ Cell = "E5"; Value = '  +1.28%  '; Quote = $false }
    @{ Cell = "D6"; Value = '166.98'; Quote = $true }
    @{ Cell = "E6"; Value = '  +0.38%  '; Quote = $false }
    @{ Cell = "E7"; Value = '  -0.11%  '; Quote = $false }
    @{ Cell = "E8"; Value = '  +0.10%  '; Quote = $false }
    @{ Cell = "E9"; Value = '  +0.58%  '; Quote = $false }
    @{ Cell = "E10"; Value = '  -0.84%  '; Quote = $false }
    @{ Cell = "E11"; Value = '  +0.86%  '; Quote = $false }
    @{ Cell = "E12"; Value = '  -0.56%  '; Quote = $false }
    @{ Cell = "D13"; Value = '35.96'; Quote = $true }
    @{ Cell = "E13"; Value = '  -0.80%  '; Quote = $false }
    @{ Cell = "D14"; Value = '4.469.48'; Quote = $false }
    @{ Cell = "E14"; Value = '  +1.18%  '; Quote = $false }
    @{ Cell = "D15"; Value = '3.859.57'; Quote = $false }
    @{ Cell = "E15"; Value = '  +2.14%  '; Quote = $false }
    @{ Cell = "E16"; Value = '  +0.35%  '; Quote = $false }
    @{ Cell = "D17"; Value = '67.848.90'; Quote = $false }
    @{ Cell = "E17"; Value = '  +0.37%  '; Quote = $false }
    @{ Cell = "E18"; Value = '  +1.42%  '; Quote = $false }
    @{ Cell = "D20"; Value = '463.58'; Quote = $true }
    @{ Cell = "E20"; Value = '  +1.42%  '; Quote = $false }
    @{ Cell = "D21"; Value = '9.95'; Quote = $true }
    @{ Cell = "E21"; Value = '  -1.33%  '; Quote = $false }
    @{ Cell = "D22"; Value = '0.702'; Quote = $true }
    @{ Cell = "E22"; Value = '  +0.54%  '; Quote = $false }
    @{ Cell = "E23"; Value = '  -3.42%  '; Quote = $false }
    @{ Cell = "D24"; Value = '83.46'; Quote = $true }
    @{ Cell = "E24"; Value = '  +0.13%  '; Quote = $false }
    @{ Cell = "E25"; Value = '  +1.35%  '; Quote = $false }
    @{ Cell = "E26"; Value = '  -0.96%  '; Quote = $false }
    @{ Cell = "D27"; Value = '10.11'; Quote = $true }
    @{ Cell = "E27"; Value = '  +0.34%  '; Quote = $false }
    @{ Cell = "E28"; Value = '  +0.17%  '; Quote = $false }
    @{ Cell = "D29"; Value = '3.974.54'; Quote = $false }
    @{ Cell = "E29"; Value = '  +1.16%  '; Quote = $false }
    @{ Cell = "E30"; Value = '  +0.15%  '; Quote = $false }
    @{ Cell = "D31"; Value = '7.41'; Quote = $true }
    @{ Cell = "E31"; Value = '  +1.73%  '; Quote = $false }
    @{ Cell = "E32"; Value = '  +1.80%  '; Quote = $false }
    @{ Cell = "D33"; Value = '29.75'; Quote = $true }
    @{ Cell = "E33"; Value = '  -0.17%  '; Quote = $false }
    @{ Cell = "D34"; Value = '0.999'; Quote = $true }
    @{ Cell = "E34"; Value = '  +0.00%  '; Quote = $false }
    @{ Cell = "E35"; Value = '  -1.21%  '; Quote = $false }
    @{ Cell = "E36"; Value = '  +0.27%  '; Quote = $false }
    @{ Cell = "D37"; Value = '3.35'; Quote = $true }
    @{ Cell = "E37"; Value = '  +0.10%  '; Quote = $false }
    @{ Cell = "E38"; Value = '  +0.31%  '; Quote = $false }
    @{ Cell = "D39"; Value = '1.00'; Quote = $true }
    @{ Cell = "E39"; Value = '  +0.61%  '; Quote = $false }
    @{ Cell = "E40"; Value = '  +0.88%  '; Quote = $false }
    @{ Cell = "D41"; Value = '1.00'; Quote = $true }
    @{ Cell = "E41"; Value = '  +0.04%  '; Quote = $false }
    @{ Cell = "D43"; Value = '48.10'; Quote = $true }
    @{ Cell = "E43"; Value = '  +2.17%  '; Quote = $false }
    @{ Cell = "E44"; Value = '  +0.58%  '; Quote = $false }
    @{ Cell = "B45"; Value = 'EnergySwap'; Quote = $false }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Quote = $false }
    @{ Cell = "D45"; Value = '28.53'; Quote = $true }
    @{ Cell = "E45"; Value = '  +11.37%  '; Quote = $false }
    @{ Cell = "B46"; Value = 'Arweave'; Quote = $false }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; Quote = $false }
    @{ Cell = "D46"; Value = '43.26'; Quote = $true }
    @{ Cell = "E46"; Value = '  -4.50%  '; Quote = $false }
    @{ Cell = "D47"; Value = '1.40'; Quote = $true }
    @{ Cell = "E47"; Value = '  +11.98%  '; Quote = $false }
    @{ Cell = "D48"; Value = '8.36'; Quote = $true }
    @{ Cell = "E48"; Value = '  +0.30%  '; Quote = $false }
    @{ Cell = "D49"; Value = '148.17'; Quote = $true }
    @{ Cell = "E49"; Value = '  -0.07%  '; Quote = $false }
    @{ Cell = "D50"; Value = '1.84'; Quote = $true }
    @{ Cell = "E50"; Value = '  +0.33%  '; Quote = $false }
    @{ Cell = "D51"; Value = '388.72'; Quote = $true }
    @{ Cell = "E51"; Value = '  -0.20%  '; Quote = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Quote) {
        $cell.Value = "'" + $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
